$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-6 and rows 7-9 were swapped as complete record blocks
# (the 3-row block starting at row 4 exchanged places with the
# 3-row block starting at row 7), columns A through AY.
$rngTop = $ws.Range("A4:AY6")
$rngBottom = $ws.Range("A7:AY9")

# Capture the full values of both blocks before writing anything back.
$valsTop = $rngTop.Value2
$valsBottom = $rngBottom.Value2

# Columns Y, Z, AA, AB hold date/time values stored as plain text
# (e.g. "2021-10-26", "00:00"). Force those ranges to Text format
# before assigning so Excel does not reinterpret them as date serials.
$ws.Range("Y4:AB6").NumberFormat = "@"
$ws.Range("Y7:AB9").NumberFormat = "@"

$rngTop.Value2 = $valsBottom
$rngBottom.Value2 = $valsTop
